# Added original format (hardcoded)
# Adds a new "OriginalFormat" column to the "Title" sheet, with a single
# hardcoded data value of "35mm" for the existing row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Title")

# The sheet currently has 14 columns (A..N); the new column goes in O.
$headerCell = $ws.Cells.Item(1, 15)
$dataCell   = $ws.Cells.Item(2, 15)

# New header, styled like the rest of the header row (green fill).
$headerCell.Value = "OriginalFormat"
$headerCell.Interior.Color = $ws.Cells.Item(1, 1).Interior.Color

# Hardcoded data value for the one existing title.
$dataCell.Value = "35mm"

# Match the width of the other "narrow" columns (e.g. TitleCode in A).
$ws.Range("O1").ColumnWidth = $ws.Cells.Item(1, 1).ColumnWidth
